# Scheduled runner update: refresh market-price-derived profit columns (H:N)
# across the Leve profit sheets, per latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 3460.5454
$ws.Cells.Item(17, 10).Value = 3678.5806
$ws.Cells.Item(17, 12).Value = 11035.7418
$ws.Cells.Item(17, 14).Value = -11371.7418

$ws.Cells.Item(74, 8).Value = 5228.5386
$ws.Cells.Item(74, 9).Value = 4330.3335
$ws.Cells.Item(74, 11).Value = 4330.3335
$ws.Cells.Item(74, 13).Value = -3394.3335

$ws.Cells.Item(77, 8).Value = 5228.5386
$ws.Cells.Item(77, 9).Value = 4330.3335
$ws.Cells.Item(77, 11).Value = 21651.6675
$ws.Cells.Item(77, 13).Value = -16971.6675

$ws.Cells.Item(112, 8).Value = 1470.4117

$ws.Cells.Item(129, 8).Value = 995.38776
$ws.Cells.Item(129, 9).Value = 341
$ws.Cells.Item(129, 10).Value = 1053.5555
$ws.Cells.Item(129, 11).Value = 1023
$ws.Cells.Item(129, 12).Value = 3160.6665
$ws.Cells.Item(129, 13).Value = 3977
$ws.Cells.Item(129, 14).Value = -13160.6665

$ws.Cells.Item(133, 8).Value = 63695
$ws.Cells.Item(133, 10).Value = 63695
$ws.Cells.Item(133, 12).Value = 63695
$ws.Cells.Item(133, 14).Value = -73815

$ws.Cells.Item(137, 8).Value = 2834.5454
$ws.Cells.Item(137, 9).Value = 2882.8572
$ws.Cells.Item(137, 11).Value = 8648.571599999999
$ws.Cells.Item(137, 13).Value = -6098.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1773.6842
$ws.Cells.Item(45, 9).Value = 1747.0588
$ws.Cells.Item(45, 10).Value = 2000
$ws.Cells.Item(45, 11).Value = 1747.0588
$ws.Cells.Item(45, 12).Value = 2000
$ws.Cells.Item(45, 13).Value = -1370.0588
$ws.Cells.Item(45, 14).Value = -2754

$ws.Cells.Item(97, 8).Value = 881.38464
$ws.Cells.Item(97, 9).Value = 700.7619
$ws.Cells.Item(97, 11).Value = 700.7619
$ws.Cells.Item(97, 13).Value = -204.7619

$ws.Cells.Item(110, 8).Value = 1661.0834
$ws.Cells.Item(110, 9).Value = 1492.5555
$ws.Cells.Item(110, 10).Value = 2166.6667
$ws.Cells.Item(110, 11).Value = 1492.5555
$ws.Cells.Item(110, 12).Value = 2166.6667
$ws.Cells.Item(110, 13).Value = 552.4445000000001
$ws.Cells.Item(110, 14).Value = -6256.6667

$ws.Cells.Item(122, 8).Value = 6251831.5
$ws.Cells.Item(122, 9).Value = 1841.5834
$ws.Cells.Item(122, 10).Value = 15626816
$ws.Cells.Item(122, 11).Value = 5524.7502
$ws.Cells.Item(122, 12).Value = 46880448
$ws.Cells.Item(122, 13).Value = -3074.7502
$ws.Cells.Item(122, 14).Value = -46885348

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 54593.332
$ws.Cells.Item(52, 10).Value = 54593.332
$ws.Cells.Item(52, 12).Value = 54593.332
$ws.Cells.Item(52, 14).Value = -55181.332

$ws.Cells.Item(58, 8).Value = 2221760.8
$ws.Cells.Item(58, 9).Value = 3638170.8
$ws.Cells.Item(58, 10).Value = 8620.375
$ws.Cells.Item(58, 11).Value = 3638170.8
$ws.Cells.Item(58, 12).Value = 8620.375
$ws.Cells.Item(58, 13).Value = -3637967.8
$ws.Cells.Item(58, 14).Value = -9026.375

$ws.Cells.Item(122, 8).Value = 11169.1875
$ws.Cells.Item(122, 9).Value = 6727.727
$ws.Cells.Item(122, 10).Value = 20940.4
$ws.Cells.Item(122, 11).Value = 20183.181
$ws.Cells.Item(122, 12).Value = 62821.2
$ws.Cells.Item(122, 13).Value = -17733.181
$ws.Cells.Item(122, 14).Value = -67721.20000000001

$ws.Cells.Item(136, 8).Value = 2221760.8
$ws.Cells.Item(136, 9).Value = 3638170.8
$ws.Cells.Item(136, 10).Value = 8620.375
$ws.Cells.Item(136, 11).Value = 10914512.4
$ws.Cells.Item(136, 12).Value = 25861.125
$ws.Cells.Item(136, 13).Value = -10911962.4
$ws.Cells.Item(136, 14).Value = -30961.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 2457
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 2457
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 7371
$ws.Cells.Item(62, 14).Value = -8743
$ws.Cells.Item(62, 13).ClearContents()

$ws.Cells.Item(63, 8).Value = 3968.2307
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 3968.2307
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 11904.6921
$ws.Cells.Item(63, 14).Value = -13402.6921
$ws.Cells.Item(63, 13).ClearContents()

$ws.Cells.Item(65, 8).Value = 2457
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 2457
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 22113
$ws.Cells.Item(65, 14).Value = -28977
$ws.Cells.Item(65, 13).ClearContents()

$ws.Cells.Item(66, 8).Value = 3968.2307
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 3968.2307
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 35714.0763
$ws.Cells.Item(66, 14).Value = -43202.0763
$ws.Cells.Item(66, 13).ClearContents()

$ws.Cells.Item(69, 8).Value = 71433140
$ws.Cells.Item(69, 9).Value = 2000
$ws.Cells.Item(69, 10).Value = 100005590
$ws.Cells.Item(69, 11).Value = 6000
$ws.Cells.Item(69, 12).Value = 300016770
$ws.Cells.Item(69, 13).Value = -5189
$ws.Cells.Item(69, 14).Value = -300018392

$ws.Cells.Item(70, 8).Value = 3555.182
$ws.Cells.Item(70, 9).Value = 1556
$ws.Cells.Item(70, 10).Value = 3999.4443
$ws.Cells.Item(70, 11).Value = 4668
$ws.Cells.Item(70, 12).Value = 11998.3329
$ws.Cells.Item(70, 13).Value = -4353
$ws.Cells.Item(70, 14).Value = -12628.3329

$ws.Cells.Item(72, 8).Value = 71433140
$ws.Cells.Item(72, 9).Value = 2000
$ws.Cells.Item(72, 10).Value = 100005590
$ws.Cells.Item(72, 11).Value = 18000
$ws.Cells.Item(72, 12).Value = 900050310
$ws.Cells.Item(72, 13).Value = -13944
$ws.Cells.Item(72, 14).Value = -900058422

$ws.Cells.Item(73, 8).Value = 3555.182
$ws.Cells.Item(73, 9).Value = 1556
$ws.Cells.Item(73, 10).Value = 3999.4443
$ws.Cells.Item(73, 11).Value = 4668
$ws.Cells.Item(73, 12).Value = 11998.3329
$ws.Cells.Item(73, 13).Value = -3576
$ws.Cells.Item(73, 14).Value = -14182.3329

$ws.Cells.Item(107, 8).Value = 1435.1818
$ws.Cells.Item(107, 10).Value = 1946.2667
$ws.Cells.Item(107, 12).Value = 5838.800099999999
$ws.Cells.Item(107, 14).Value = -9678.8001

$ws.Cells.Item(113, 8).Value = 694.93616
$ws.Cells.Item(113, 9).Value = 702
$ws.Cells.Item(113, 10).Value = 660.5
$ws.Cells.Item(113, 11).Value = 2106
$ws.Cells.Item(113, 12).Value = 1981.5
$ws.Cells.Item(113, 13).Value = 64
$ws.Cells.Item(113, 14).Value = -6321.5

$ws.Cells.Item(122, 8).Value = 863.63635
$ws.Cells.Item(122, 9).Value = 690
$ws.Cells.Item(122, 10).Value = 881
$ws.Cells.Item(122, 11).Value = 6210
$ws.Cells.Item(122, 12).Value = 7929
$ws.Cells.Item(122, 13).Value = -3760
$ws.Cells.Item(122, 14).Value = -12829

$ws.Cells.Item(124, 8).Value = 3012.4
$ws.Cells.Item(124, 10).Value = 3680.5715
$ws.Cells.Item(124, 12).Value = 11041.7145
$ws.Cells.Item(124, 14).Value = -20861.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2677.4614
$ws.Cells.Item(126, 9).Value = 1974.2667
$ws.Cells.Item(126, 10).Value = 3636.3635
$ws.Cells.Item(126, 11).Value = 5922.800099999999
$ws.Cells.Item(126, 12).Value = 10909.0905
$ws.Cells.Item(126, 13).Value = -3452.800099999999
$ws.Cells.Item(126, 14).Value = -15849.0905

$ws.Cells.Item(137, 8).Value = 49800
$ws.Cells.Item(137, 10).Value = 49800
$ws.Cells.Item(137, 12).Value = 49800
$ws.Cells.Item(137, 14).Value = -60000

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1620
$ws.Cells.Item(82, 9).Value = 1500
$ws.Cells.Item(82, 10).Value = 1650
$ws.Cells.Item(82, 11).Value = 1500
$ws.Cells.Item(82, 12).Value = 1650
$ws.Cells.Item(82, 13).Value = -1139
$ws.Cells.Item(82, 14).Value = -2372

$ws.Cells.Item(85, 8).Value = 1620
$ws.Cells.Item(85, 9).Value = 1500
$ws.Cells.Item(85, 10).Value = 1650
$ws.Cells.Item(85, 11).Value = 1500
$ws.Cells.Item(85, 12).Value = 1650
$ws.Cells.Item(85, 13).Value = -252
$ws.Cells.Item(85, 14).Value = -4146

$ws.Cells.Item(136, 8).Value = 3868.4546
$ws.Cells.Item(136, 9).Value = 2148.6667
$ws.Cells.Item(136, 10).Value = 6448.136
$ws.Cells.Item(136, 11).Value = 6446.000100000001
$ws.Cells.Item(136, 12).Value = 19344.408
$ws.Cells.Item(136, 13).Value = -3896.000100000001
$ws.Cells.Item(136, 14).Value = -24444.408

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1192.5294
$ws.Cells.Item(126, 9).Value = 1177.1786
$ws.Cells.Item(126, 10).Value = 1264.1666
$ws.Cells.Item(126, 11).Value = 3531.5358
$ws.Cells.Item(126, 12).Value = 3792.4998
$ws.Cells.Item(126, 13).Value = -1061.5358
$ws.Cells.Item(126, 14).Value = -8732.4998

$ws.Cells.Item(135, 8).Value = 45998.75
$ws.Cells.Item(135, 10).Value = 45998.75
$ws.Cells.Item(135, 12).Value = 45998.75
$ws.Cells.Item(135, 14).Value = -56138.75
